$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9565056562423706
$ws.Range("B1").Value = 3.025663137435913
$ws.Range("C1").Value = 4.198607921600342
$ws.Range("D1").Value = 2.073939561843872
$ws.Range("E1").Value = 1.235293984413147
